$d = $word.ActiveDocument

function Check-Found($found, $label) {
    if (-not $found) {
        throw "Find failed: $label"
    }
}

# --- Change 1 -------------------------------------------------------------
# "Svakim restoranom mora upravljati jedan ili vi[še menadžer]a. "
#   -> "Svakim restoranom može upravljati vi[še menadžer]a ili ni jedan. "
# Replace " mora upravljati jedan ili vi" (two plain, non-bold runs) with
# " može upravljati vi" (stays plain / non-bold), directly after the bold
# "restoranom" run.
$r1 = $d.Content
$found1 = $r1.Find.Execute(" mora upravljati jedan ili vi", $true, $false, $false, $false, $false, $true, 1, $false, " može upravljati vi", 2)
Check-Found $found1 "mora upravljati jedan ili vi"

# --- Change 2 ---------------------------------------------------------
# Append " ili ni jedan" right after "...menadžera" (before the trailing
# ". "). Search only after the text touched by change 1 so we land on the
# right "menadžera" (the word also occurs later, unrelated, in the doc).
$r2 = $d.Range($r1.End, $d.Content.End)
$found2 = $r2.Find.Execute("menadžera", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Check-Found $found2 "menadžera"
$r2.Collapse(0)
$r2.InsertAfter(" ili ni jedan")

# The leading space inherits the bold formatting already present at the
# insertion point (end of the bold "menadžera"), matching the diff's bold
# trailing space run; the following words are then set back to non-bold.
$r3 = $d.Range($r2.Start, $d.Content.End)
$found3 = $r3.Find.Execute("ili ni jedan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Check-Found $found3 "ili ni jedan"
$r3.Font.Bold = $false

# --- Change 3 -------------------------------------------------------------
# "Svaki restoran mora biti čuvan od strane najmanje jednog portira. "
#   -> "Svaki restoran može biti čuvan od strane više portira ili da ne
#       bude čuvan uopšte. "
# Leave the bold "restoran" run untouched; replace everything from the
# space after it up to (and including) "portira" with new plain text.
$r4 = $d.Range($r3.End, $d.Content.End)
$found4 = $r4.Find.Execute(" mora biti čuvan od strane najmanje jednog portira", $true, $false, $false, $false, $false, $true, 1, $false, " može biti čuvan od strane više portira ili da ne bude čuvan uopšte", 2)
Check-Found $found4 "mora biti cuvan ... portira"

Write-Host "All replacements applied: found1=$found1 found2=$found2 found3=$found3 found4=$found4"
